# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N, which shifts the existing "Late" / "heading" / "Outstanding"
# columns one place to the right and adds a new blank header/column in
# their place. The sheet also becomes the active tab/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (moves tabSelected/activeTab here).
$ws.Activate()

# Insert a new column before column N (14th column).
$ws.Columns.Item(14).Insert()

# The newly inserted column inherits the width of its left neighbour
# (column M), which was 11 characters wide.
$ws.Columns.Item(14).ColumnWidth = 10.2

# Leave the selection where it ended up after the edit.
[void]$ws.Range("R6").Select()
